# ============================================================
# 2022-Q1 fund-holdings update.
#
# The workbook keeps one sheet per quarter plus a rolling "总计"
# (grand totals) sheet at the end. To add a new quarter we:
#   1. Repurpose the current (last) sheet -- today that is "总计" --
#      as the detail sheet for the new quarter, "2022-Q1".
#   2. Append a brand-new "总计" sheet after it, re-populated with
#      the prior totals history plus a new leading row for 2022-Q1.
# This mirrors exactly how the previous quarters were appended, and
# keeps the grand-totals sheet last in tab order.
# ============================================================
$wb = $excel.ActiveWorkbook

# ---- Step 1: turn the old "总计" sheet into the "2022-Q1" detail sheet ----
$detailSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$detailSheet.Name = "2022-Q1"
$detailSheet.Cells.Clear()

# Style donors: an existing quarterly sheet already has the exact header
# and index-column formatting we need, so copy styles from there instead
# of re-building them by hand.
$styleSheet = $wb.Worksheets.Item(1)

# ---- Header row (bold, centered, bordered) ----
$styleSheet.Range("B1").Copy($detailSheet.Range("B1"))
$detailSheet.Range("B1").Value = "基金代码"
$styleSheet.Range("B1").Copy($detailSheet.Range("C1"))
$detailSheet.Range("C1").Value = "基金名称"
$styleSheet.Range("B1").Copy($detailSheet.Range("D1"))
$detailSheet.Range("D1").Value = "基金规模"
$styleSheet.Range("B1").Copy($detailSheet.Range("E1"))
$detailSheet.Range("E1").Value = "股票总仓位"
$styleSheet.Range("B1").Copy($detailSheet.Range("F1"))
$detailSheet.Range("F1").Value = "仓位占比"
$styleSheet.Range("B1").Copy($detailSheet.Range("G1"))
$detailSheet.Range("G1").Value = "持有市值(亿元)"
$styleSheet.Range("B1").Copy($detailSheet.Range("H1"))
$detailSheet.Range("H1").Value = "仓位排名"

# ---- Data rows ----
# Row 2
$styleSheet.Range("A2").Copy($detailSheet.Range("A2"))
$detailSheet.Range("A2").Value = 0
$detailSheet.Range("B2").NumberFormat = "@"
$detailSheet.Range("B2").Value = "011174"
$detailSheet.Range("B2").Style = "Normal"
$detailSheet.Range("D2").NumberFormat = "@"
$detailSheet.Range("D2").Value = "63.01"
$detailSheet.Range("D2").Style = "Normal"
$detailSheet.Range("E2").NumberFormat = "@"
$detailSheet.Range("E2").Value = "93.47"
$detailSheet.Range("E2").Style = "Normal"
$detailSheet.Range("F2").NumberFormat = "@"
$detailSheet.Range("F2").Value = "8.14"
$detailSheet.Range("F2").Style = "Normal"
$detailSheet.Range("G2").NumberFormat = "@"
$detailSheet.Range("G2").Value = "5.1290"
$detailSheet.Range("G2").Style = "Normal"
$detailSheet.Range("C2").Value = "中庚价值品质一年持有期混合"
$detailSheet.Range("H2").Value = 3

# Row 3
$styleSheet.Range("A2").Copy($detailSheet.Range("A3"))
$detailSheet.Range("A3").Value = 1
$detailSheet.Range("B3").NumberFormat = "@"
$detailSheet.Range("B3").Value = "007130"
$detailSheet.Range("B3").Style = "Normal"
$detailSheet.Range("D3").NumberFormat = "@"
$detailSheet.Range("D3").Value = "40.99"
$detailSheet.Range("D3").Style = "Normal"
$detailSheet.Range("E3").NumberFormat = "@"
$detailSheet.Range("E3").Value = "93.10"
$detailSheet.Range("E3").Style = "Normal"
$detailSheet.Range("F3").NumberFormat = "@"
$detailSheet.Range("F3").Value = "7.12"
$detailSheet.Range("F3").Style = "Normal"
$detailSheet.Range("G3").NumberFormat = "@"
$detailSheet.Range("G3").Value = "2.9185"
$detailSheet.Range("G3").Style = "Normal"
$detailSheet.Range("C3").Value = "中庚小盘价值股票"
$detailSheet.Range("H3").Value = 2

# Row 4
$styleSheet.Range("A2").Copy($detailSheet.Range("A4"))
$detailSheet.Range("A4").Value = 2
$detailSheet.Range("B4").NumberFormat = "@"
$detailSheet.Range("B4").Value = "006551"
$detailSheet.Range("B4").Style = "Normal"
$detailSheet.Range("D4").NumberFormat = "@"
$detailSheet.Range("D4").Value = "36.49"
$detailSheet.Range("D4").Style = "Normal"
$detailSheet.Range("E4").NumberFormat = "@"
$detailSheet.Range("E4").Value = "94.18"
$detailSheet.Range("E4").Style = "Normal"
$detailSheet.Range("F4").NumberFormat = "@"
$detailSheet.Range("F4").Value = "6.33"
$detailSheet.Range("F4").Style = "Normal"
$detailSheet.Range("G4").NumberFormat = "@"
$detailSheet.Range("G4").Value = "2.3098"
$detailSheet.Range("G4").Style = "Normal"
$detailSheet.Range("C4").Value = "中庚价值领航混合"
$detailSheet.Range("H4").Value = 5

# Row 5
$styleSheet.Range("A2").Copy($detailSheet.Range("A5"))
$detailSheet.Range("A5").Value = 3
$detailSheet.Range("B5").NumberFormat = "@"
$detailSheet.Range("B5").Value = "007497"
$detailSheet.Range("B5").Style = "Normal"
$detailSheet.Range("D5").NumberFormat = "@"
$detailSheet.Range("D5").Value = "24.35"
$detailSheet.Range("D5").Style = "Normal"
$detailSheet.Range("E5").NumberFormat = "@"
$detailSheet.Range("E5").Value = "89.42"
$detailSheet.Range("E5").Style = "Normal"
$detailSheet.Range("F5").NumberFormat = "@"
$detailSheet.Range("F5").Value = "5.00"
$detailSheet.Range("F5").Style = "Normal"
$detailSheet.Range("G5").NumberFormat = "@"
$detailSheet.Range("G5").Value = "1.2175"
$detailSheet.Range("G5").Style = "Normal"
$detailSheet.Range("C5").Value = "中庚价值灵动灵活配置混合"
$detailSheet.Range("H5").Value = 4

# Row 6
$styleSheet.Range("A2").Copy($detailSheet.Range("A6"))
$detailSheet.Range("A6").Value = 4
$detailSheet.Range("B6").NumberFormat = "@"
$detailSheet.Range("B6").Value = "410003"
$detailSheet.Range("B6").Style = "Normal"
$detailSheet.Range("D6").NumberFormat = "@"
$detailSheet.Range("D6").Value = "11.98"
$detailSheet.Range("D6").Style = "Normal"
$detailSheet.Range("E6").NumberFormat = "@"
$detailSheet.Range("E6").Value = "94.90"
$detailSheet.Range("E6").Style = "Normal"
$detailSheet.Range("F6").NumberFormat = "@"
$detailSheet.Range("F6").Value = "4.76"
$detailSheet.Range("F6").Style = "Normal"
$detailSheet.Range("G6").NumberFormat = "@"
$detailSheet.Range("G6").Value = "0.5702"
$detailSheet.Range("G6").Style = "Normal"
$detailSheet.Range("C6").Value = "华富成长趋势混合"
$detailSheet.Range("H6").Value = 7

# Row 7
$styleSheet.Range("A2").Copy($detailSheet.Range("A7"))
$detailSheet.Range("A7").Value = 5
$detailSheet.Range("B7").NumberFormat = "@"
$detailSheet.Range("B7").Value = "410007"
$detailSheet.Range("B7").Style = "Normal"
$detailSheet.Range("D7").NumberFormat = "@"
$detailSheet.Range("D7").Value = "11.15"
$detailSheet.Range("D7").Style = "Normal"
$detailSheet.Range("E7").NumberFormat = "@"
$detailSheet.Range("E7").Value = "79.40"
$detailSheet.Range("E7").Style = "Normal"
$detailSheet.Range("F7").NumberFormat = "@"
$detailSheet.Range("F7").Value = "4.64"
$detailSheet.Range("F7").Style = "Normal"
$detailSheet.Range("G7").NumberFormat = "@"
$detailSheet.Range("G7").Value = "0.5174"
$detailSheet.Range("G7").Style = "Normal"
$detailSheet.Range("C7").Value = "华富价值增长混合"
$detailSheet.Range("H7").Value = 5

# Row 8
$styleSheet.Range("A2").Copy($detailSheet.Range("A8"))
$detailSheet.Range("A8").Value = 6
$detailSheet.Range("B8").NumberFormat = "@"
$detailSheet.Range("B8").Value = "014024"
$detailSheet.Range("B8").Style = "Normal"
$detailSheet.Range("D8").NumberFormat = "@"
$detailSheet.Range("D8").Value = "8.27"
$detailSheet.Range("D8").Style = "Normal"
$detailSheet.Range("E8").NumberFormat = "@"
$detailSheet.Range("E8").Value = "89.06"
$detailSheet.Range("E8").Style = "Normal"
$detailSheet.Range("F8").NumberFormat = "@"
$detailSheet.Range("F8").Value = "4.33"
$detailSheet.Range("F8").Style = "Normal"
$detailSheet.Range("G8").NumberFormat = "@"
$detailSheet.Range("G8").Value = "0.3581"
$detailSheet.Range("G8").Style = "Normal"
$detailSheet.Range("C8").Value = "华富卓越成长一年持有期混合A"
$detailSheet.Range("H8").Value = 8

# Row 9
$styleSheet.Range("A2").Copy($detailSheet.Range("A9"))
$detailSheet.Range("A9").Value = 7
$detailSheet.Range("B9").NumberFormat = "@"
$detailSheet.Range("B9").Value = "410001"
$detailSheet.Range("B9").Style = "Normal"
$detailSheet.Range("D9").NumberFormat = "@"
$detailSheet.Range("D9").Value = "3.34"
$detailSheet.Range("D9").Style = "Normal"
$detailSheet.Range("E9").NumberFormat = "@"
$detailSheet.Range("E9").Value = "89.23"
$detailSheet.Range("E9").Style = "Normal"
$detailSheet.Range("F9").NumberFormat = "@"
$detailSheet.Range("F9").Value = "9.12"
$detailSheet.Range("F9").Style = "Normal"
$detailSheet.Range("G9").NumberFormat = "@"
$detailSheet.Range("G9").Value = "0.3046"
$detailSheet.Range("G9").Style = "Normal"
$detailSheet.Range("C9").Value = "华富竞争力优选混合"
$detailSheet.Range("H9").Value = 1

# Row 10
$styleSheet.Range("A2").Copy($detailSheet.Range("A10"))
$detailSheet.Range("A10").Value = 8
$detailSheet.Range("B10").NumberFormat = "@"
$detailSheet.Range("B10").Value = "515150"
$detailSheet.Range("B10").Style = "Normal"
$detailSheet.Range("D10").NumberFormat = "@"
$detailSheet.Range("D10").Value = "7.41"
$detailSheet.Range("D10").Style = "Normal"
$detailSheet.Range("E10").NumberFormat = "@"
$detailSheet.Range("E10").Value = "99.21"
$detailSheet.Range("E10").Style = "Normal"
$detailSheet.Range("F10").NumberFormat = "@"
$detailSheet.Range("F10").Value = "3.57"
$detailSheet.Range("F10").Style = "Normal"
$detailSheet.Range("G10").NumberFormat = "@"
$detailSheet.Range("G10").Value = "0.2645"
$detailSheet.Range("G10").Style = "Normal"
$detailSheet.Range("C10").Value = "富国中证国企一带一路ETF"
$detailSheet.Range("H10").Value = 1

# Row 11
$styleSheet.Range("A2").Copy($detailSheet.Range("A11"))
$detailSheet.Range("A11").Value = 9
$detailSheet.Range("B11").NumberFormat = "@"
$detailSheet.Range("B11").Value = "009398"
$detailSheet.Range("B11").Style = "Normal"
$detailSheet.Range("D11").NumberFormat = "@"
$detailSheet.Range("D11").Value = "5.40"
$detailSheet.Range("D11").Style = "Normal"
$detailSheet.Range("E11").NumberFormat = "@"
$detailSheet.Range("E11").Value = "93.67"
$detailSheet.Range("E11").Style = "Normal"
$detailSheet.Range("F11").NumberFormat = "@"
$detailSheet.Range("F11").Value = "4.61"
$detailSheet.Range("F11").Style = "Normal"
$detailSheet.Range("G11").NumberFormat = "@"
$detailSheet.Range("G11").Value = "0.2489"
$detailSheet.Range("G11").Style = "Normal"
$detailSheet.Range("C11").Value = "华富成长企业精选股票"
$detailSheet.Range("H11").Value = 8

# Row 12
$styleSheet.Range("A2").Copy($detailSheet.Range("A12"))
$detailSheet.Range("A12").Value = 10
$detailSheet.Range("B12").NumberFormat = "@"
$detailSheet.Range("B12").Value = "002076"
$detailSheet.Range("B12").Style = "Normal"
$detailSheet.Range("D12").NumberFormat = "@"
$detailSheet.Range("D12").Value = "14.53"
$detailSheet.Range("D12").Style = "Normal"
$detailSheet.Range("E12").NumberFormat = "@"
$detailSheet.Range("E12").Value = "93.68"
$detailSheet.Range("E12").Style = "Normal"
$detailSheet.Range("F12").NumberFormat = "@"
$detailSheet.Range("F12").Value = "1.62"
$detailSheet.Range("F12").Style = "Normal"
$detailSheet.Range("G12").NumberFormat = "@"
$detailSheet.Range("G12").Value = "0.2354"
$detailSheet.Range("G12").Style = "Normal"
$detailSheet.Range("C12").Value = "浙商中证500指数增强A"
$detailSheet.Range("H12").Value = 3

# Row 13
$styleSheet.Range("A2").Copy($detailSheet.Range("A13"))
$detailSheet.Range("A13").Value = 11
$detailSheet.Range("B13").NumberFormat = "@"
$detailSheet.Range("B13").Value = "515110"
$detailSheet.Range("B13").Style = "Normal"
$detailSheet.Range("D13").NumberFormat = "@"
$detailSheet.Range("D13").Value = "4.83"
$detailSheet.Range("D13").Style = "Normal"
$detailSheet.Range("E13").NumberFormat = "@"
$detailSheet.Range("E13").Value = "99.52"
$detailSheet.Range("E13").Style = "Normal"
$detailSheet.Range("F13").NumberFormat = "@"
$detailSheet.Range("F13").Value = "3.57"
$detailSheet.Range("F13").Style = "Normal"
$detailSheet.Range("G13").NumberFormat = "@"
$detailSheet.Range("G13").Value = "0.1724"
$detailSheet.Range("G13").Style = "Normal"
$detailSheet.Range("C13").Value = "易方达中证国企一带一路ETF"
$detailSheet.Range("H13").Value = 1

# Row 14
$styleSheet.Range("A2").Copy($detailSheet.Range("A14"))
$detailSheet.Range("A14").Value = 12
$detailSheet.Range("B14").NumberFormat = "@"
$detailSheet.Range("B14").Value = "460009"
$detailSheet.Range("B14").Style = "Normal"
$detailSheet.Range("D14").NumberFormat = "@"
$detailSheet.Range("D14").Value = "9.13"
$detailSheet.Range("D14").Style = "Normal"
$detailSheet.Range("E14").NumberFormat = "@"
$detailSheet.Range("E14").Value = "90.47"
$detailSheet.Range("E14").Style = "Normal"
$detailSheet.Range("F14").NumberFormat = "@"
$detailSheet.Range("F14").Value = "1.03"
$detailSheet.Range("F14").Style = "Normal"
$detailSheet.Range("G14").NumberFormat = "@"
$detailSheet.Range("G14").Value = "0.0940"
$detailSheet.Range("G14").Style = "Normal"
$detailSheet.Range("C14").Value = "华泰柏瑞量化先行混合A"
$detailSheet.Range("H14").Value = 2

# Row 15
$styleSheet.Range("A2").Copy($detailSheet.Range("A15"))
$detailSheet.Range("A15").Value = 13
$detailSheet.Range("B15").NumberFormat = "@"
$detailSheet.Range("B15").Value = "003152"
$detailSheet.Range("B15").Style = "Normal"
$detailSheet.Range("D15").NumberFormat = "@"
$detailSheet.Range("D15").Value = "1.59"
$detailSheet.Range("D15").Style = "Normal"
$detailSheet.Range("E15").NumberFormat = "@"
$detailSheet.Range("E15").Value = "93.75"
$detailSheet.Range("E15").Style = "Normal"
$detailSheet.Range("F15").NumberFormat = "@"
$detailSheet.Range("F15").Value = "4.84"
$detailSheet.Range("F15").Style = "Normal"
$detailSheet.Range("G15").NumberFormat = "@"
$detailSheet.Range("G15").Value = "0.0770"
$detailSheet.Range("G15").Style = "Normal"
$detailSheet.Range("C15").Value = "华富天鑫灵活配置混合A"
$detailSheet.Range("H15").Value = 6

# Row 16
$styleSheet.Range("A2").Copy($detailSheet.Range("A16"))
$detailSheet.Range("A16").Value = 14
$detailSheet.Range("B16").NumberFormat = "@"
$detailSheet.Range("B16").Value = "006022"
$detailSheet.Range("B16").Style = "Normal"
$detailSheet.Range("D16").NumberFormat = "@"
$detailSheet.Range("D16").Value = "3.47"
$detailSheet.Range("D16").Style = "Normal"
$detailSheet.Range("E16").NumberFormat = "@"
$detailSheet.Range("E16").Value = "86.62"
$detailSheet.Range("E16").Style = "Normal"
$detailSheet.Range("F16").NumberFormat = "@"
$detailSheet.Range("F16").Value = "1.62"
$detailSheet.Range("F16").Style = "Normal"
$detailSheet.Range("G16").NumberFormat = "@"
$detailSheet.Range("G16").Value = "0.0562"
$detailSheet.Range("G16").Style = "Normal"
$detailSheet.Range("C16").Value = "富国大盘价值量化精选混合"
$detailSheet.Range("H16").Value = 6

# Row 17
$styleSheet.Range("A2").Copy($detailSheet.Range("A17"))
$detailSheet.Range("A17").Value = 15
$detailSheet.Range("B17").NumberFormat = "@"
$detailSheet.Range("B17").Value = "007386"
$detailSheet.Range("B17").Style = "Normal"
$detailSheet.Range("D17").NumberFormat = "@"
$detailSheet.Range("D17").Value = "3.38"
$detailSheet.Range("D17").Style = "Normal"
$detailSheet.Range("E17").NumberFormat = "@"
$detailSheet.Range("E17").Value = "93.68"
$detailSheet.Range("E17").Style = "Normal"
$detailSheet.Range("F17").NumberFormat = "@"
$detailSheet.Range("F17").Value = "1.62"
$detailSheet.Range("F17").Style = "Normal"
$detailSheet.Range("G17").NumberFormat = "@"
$detailSheet.Range("G17").Value = "0.0548"
$detailSheet.Range("G17").Style = "Normal"
$detailSheet.Range("C17").Value = "浙商中证500指数增强C"
$detailSheet.Range("H17").Value = 3

# Row 18
$styleSheet.Range("A2").Copy($detailSheet.Range("A18"))
$detailSheet.Range("A18").Value = 16
$detailSheet.Range("B18").NumberFormat = "@"
$detailSheet.Range("B18").Value = "515990"
$detailSheet.Range("B18").Style = "Normal"
$detailSheet.Range("D18").NumberFormat = "@"
$detailSheet.Range("D18").Value = "1.08"
$detailSheet.Range("D18").Style = "Normal"
$detailSheet.Range("E18").NumberFormat = "@"
$detailSheet.Range("E18").Value = "99.16"
$detailSheet.Range("E18").Style = "Normal"
$detailSheet.Range("F18").NumberFormat = "@"
$detailSheet.Range("F18").Value = "3.60"
$detailSheet.Range("F18").Style = "Normal"
$detailSheet.Range("G18").NumberFormat = "@"
$detailSheet.Range("G18").Value = "0.0389"
$detailSheet.Range("G18").Style = "Normal"
$detailSheet.Range("C18").Value = "汇添富中证国企一带一路ETF"
$detailSheet.Range("H18").Value = 1

# Row 19
$styleSheet.Range("A2").Copy($detailSheet.Range("A19"))
$detailSheet.Range("A19").Value = 17
$detailSheet.Range("B19").NumberFormat = "@"
$detailSheet.Range("B19").Value = "001244"
$detailSheet.Range("B19").Style = "Normal"
$detailSheet.Range("D19").NumberFormat = "@"
$detailSheet.Range("D19").Value = "3.50"
$detailSheet.Range("D19").Style = "Normal"
$detailSheet.Range("E19").NumberFormat = "@"
$detailSheet.Range("E19").Value = "91.02"
$detailSheet.Range("E19").Style = "Normal"
$detailSheet.Range("F19").NumberFormat = "@"
$detailSheet.Range("F19").Value = "0.80"
$detailSheet.Range("F19").Style = "Normal"
$detailSheet.Range("G19").NumberFormat = "@"
$detailSheet.Range("G19").Value = "0.0280"
$detailSheet.Range("G19").Style = "Normal"
$detailSheet.Range("C19").Value = "华泰柏瑞量化智慧灵活配置混合A"
$detailSheet.Range("H19").Value = 8

# Row 20
$styleSheet.Range("A2").Copy($detailSheet.Range("A20"))
$detailSheet.Range("A20").Value = 18
$detailSheet.Range("B20").NumberFormat = "@"
$detailSheet.Range("B20").Value = "007713"
$detailSheet.Range("B20").Style = "Normal"
$detailSheet.Range("D20").NumberFormat = "@"
$detailSheet.Range("D20").Value = "0.56"
$detailSheet.Range("D20").Style = "Normal"
$detailSheet.Range("E20").NumberFormat = "@"
$detailSheet.Range("E20").Value = "86.98"
$detailSheet.Range("E20").Style = "Normal"
$detailSheet.Range("F20").NumberFormat = "@"
$detailSheet.Range("F20").Value = "4.82"
$detailSheet.Range("F20").Style = "Normal"
$detailSheet.Range("G20").NumberFormat = "@"
$detailSheet.Range("G20").Value = "0.0270"
$detailSheet.Range("G20").Style = "Normal"
$detailSheet.Range("C20").Value = "华富科技动能混合"
$detailSheet.Range("H20").Value = 8

# Row 21
$styleSheet.Range("A2").Copy($detailSheet.Range("A21"))
$detailSheet.Range("A21").Value = 19
$detailSheet.Range("B21").NumberFormat = "@"
$detailSheet.Range("B21").Value = "001917"
$detailSheet.Range("B21").Style = "Normal"
$detailSheet.Range("D21").NumberFormat = "@"
$detailSheet.Range("D21").Value = "2.33"
$detailSheet.Range("D21").Style = "Normal"
$detailSheet.Range("E21").NumberFormat = "@"
$detailSheet.Range("E21").Value = "94.20"
$detailSheet.Range("E21").Style = "Normal"
$detailSheet.Range("F21").NumberFormat = "@"
$detailSheet.Range("F21").Value = "1.14"
$detailSheet.Range("F21").Style = "Normal"
$detailSheet.Range("G21").NumberFormat = "@"
$detailSheet.Range("G21").Value = "0.0266"
$detailSheet.Range("G21").Style = "Normal"
$detailSheet.Range("C21").Value = "招商量化精选股票A"
$detailSheet.Range("H21").Value = 10

# Row 22
$styleSheet.Range("A2").Copy($detailSheet.Range("A22"))
$detailSheet.Range("A22").Value = 20
$detailSheet.Range("B22").NumberFormat = "@"
$detailSheet.Range("B22").Value = "003153"
$detailSheet.Range("B22").Style = "Normal"
$detailSheet.Range("D22").NumberFormat = "@"
$detailSheet.Range("D22").Value = "0.49"
$detailSheet.Range("D22").Style = "Normal"
$detailSheet.Range("E22").NumberFormat = "@"
$detailSheet.Range("E22").Value = "93.75"
$detailSheet.Range("E22").Style = "Normal"
$detailSheet.Range("F22").NumberFormat = "@"
$detailSheet.Range("F22").Value = "4.84"
$detailSheet.Range("F22").Style = "Normal"
$detailSheet.Range("G22").NumberFormat = "@"
$detailSheet.Range("G22").Value = "0.0237"
$detailSheet.Range("G22").Style = "Normal"
$detailSheet.Range("C22").Value = "华富天鑫灵活配置混合C"
$detailSheet.Range("H22").Value = 6

# Row 23
$styleSheet.Range("A2").Copy($detailSheet.Range("A23"))
$detailSheet.Range("A23").Value = 21
$detailSheet.Range("B23").NumberFormat = "@"
$detailSheet.Range("B23").Value = "014025"
$detailSheet.Range("B23").Style = "Normal"
$detailSheet.Range("D23").NumberFormat = "@"
$detailSheet.Range("D23").Value = "0.21"
$detailSheet.Range("D23").Style = "Normal"
$detailSheet.Range("E23").NumberFormat = "@"
$detailSheet.Range("E23").Value = "89.06"
$detailSheet.Range("E23").Style = "Normal"
$detailSheet.Range("F23").NumberFormat = "@"
$detailSheet.Range("F23").Value = "4.33"
$detailSheet.Range("F23").Style = "Normal"
$detailSheet.Range("G23").NumberFormat = "@"
$detailSheet.Range("G23").Value = "0.0091"
$detailSheet.Range("G23").Style = "Normal"
$detailSheet.Range("C23").Value = "华富卓越成长一年持有期混合C"
$detailSheet.Range("H23").Value = 8

# Row 24
$styleSheet.Range("A2").Copy($detailSheet.Range("A24"))
$detailSheet.Range("A24").Value = 22
$detailSheet.Range("B24").NumberFormat = "@"
$detailSheet.Range("B24").Value = "006104"
$detailSheet.Range("B24").Style = "Normal"
$detailSheet.Range("D24").NumberFormat = "@"
$detailSheet.Range("D24").Value = "0.84"
$detailSheet.Range("D24").Style = "Normal"
$detailSheet.Range("E24").NumberFormat = "@"
$detailSheet.Range("E24").Value = "91.02"
$detailSheet.Range("E24").Style = "Normal"
$detailSheet.Range("F24").NumberFormat = "@"
$detailSheet.Range("F24").Value = "0.80"
$detailSheet.Range("F24").Style = "Normal"
$detailSheet.Range("G24").NumberFormat = "@"
$detailSheet.Range("G24").Value = "0.0067"
$detailSheet.Range("G24").Style = "Normal"
$detailSheet.Range("C24").Value = "华泰柏瑞量化智慧灵活配置混合C"
$detailSheet.Range("H24").Value = 8

# Row 25
$styleSheet.Range("A2").Copy($detailSheet.Range("A25"))
$detailSheet.Range("A25").Value = 23
$detailSheet.Range("B25").NumberFormat = "@"
$detailSheet.Range("B25").Value = "007950"
$detailSheet.Range("B25").Style = "Normal"
$detailSheet.Range("D25").NumberFormat = "@"
$detailSheet.Range("D25").Value = "0.56"
$detailSheet.Range("D25").Style = "Normal"
$detailSheet.Range("E25").NumberFormat = "@"
$detailSheet.Range("E25").Value = "94.20"
$detailSheet.Range("E25").Style = "Normal"
$detailSheet.Range("F25").NumberFormat = "@"
$detailSheet.Range("F25").Value = "1.14"
$detailSheet.Range("F25").Style = "Normal"
$detailSheet.Range("G25").NumberFormat = "@"
$detailSheet.Range("G25").Value = "0.0064"
$detailSheet.Range("G25").Style = "Normal"
$detailSheet.Range("C25").Value = "招商量化精选股票C"
$detailSheet.Range("H25").Value = 10

# Row 26
$styleSheet.Range("A2").Copy($detailSheet.Range("A26"))
$detailSheet.Range("A26").Value = 24
$detailSheet.Range("B26").NumberFormat = "@"
$detailSheet.Range("B26").Value = "005616"
$detailSheet.Range("B26").Style = "Normal"
$detailSheet.Range("D26").NumberFormat = "@"
$detailSheet.Range("D26").Value = "0.20"
$detailSheet.Range("D26").Style = "Normal"
$detailSheet.Range("E26").NumberFormat = "@"
$detailSheet.Range("E26").Value = "92.24"
$detailSheet.Range("E26").Style = "Normal"
$detailSheet.Range("F26").NumberFormat = "@"
$detailSheet.Range("F26").Value = "2.65"
$detailSheet.Range("F26").Style = "Normal"
$detailSheet.Range("G26").NumberFormat = "@"
$detailSheet.Range("G26").Value = "0.0053"
$detailSheet.Range("G26").Style = "Normal"
$detailSheet.Range("C26").Value = "东方量化成长灵活配置混合"
$detailSheet.Range("H26").Value = 2

# Row 27
$styleSheet.Range("A2").Copy($detailSheet.Range("A27"))
$detailSheet.Range("A27").Value = 25
$detailSheet.Range("B27").NumberFormat = "@"
$detailSheet.Range("B27").Value = "515510"
$detailSheet.Range("B27").Style = "Normal"
$detailSheet.Range("D27").NumberFormat = "@"
$detailSheet.Range("D27").Value = "0.15"
$detailSheet.Range("D27").Style = "Normal"
$detailSheet.Range("E27").NumberFormat = "@"
$detailSheet.Range("E27").Value = "98.79"
$detailSheet.Range("E27").Style = "Normal"
$detailSheet.Range("F27").NumberFormat = "@"
$detailSheet.Range("F27").Value = "1.22"
$detailSheet.Range("F27").Style = "Normal"
$detailSheet.Range("G27").NumberFormat = "@"
$detailSheet.Range("G27").Value = "0.0018"
$detailSheet.Range("G27").Style = "Normal"
$detailSheet.Range("C27").Value = "嘉实中证500成长估值ETF"
$detailSheet.Range("H27").Value = 8

# Row 28
$styleSheet.Range("A2").Copy($detailSheet.Range("A28"))
$detailSheet.Range("A28").Value = 26
$detailSheet.Range("B28").NumberFormat = "@"
$detailSheet.Range("B28").Value = "010246"
$detailSheet.Range("B28").Style = "Normal"
$detailSheet.Range("D28").NumberFormat = "@"
$detailSheet.Range("D28").Value = "0.12"
$detailSheet.Range("D28").Style = "Normal"
$detailSheet.Range("E28").NumberFormat = "@"
$detailSheet.Range("E28").Value = "90.47"
$detailSheet.Range("E28").Style = "Normal"
$detailSheet.Range("F28").NumberFormat = "@"
$detailSheet.Range("F28").Value = "1.03"
$detailSheet.Range("F28").Style = "Normal"
$detailSheet.Range("G28").NumberFormat = "@"
$detailSheet.Range("G28").Value = "0.0012"
$detailSheet.Range("G28").Style = "Normal"
$detailSheet.Range("C28").Value = "华泰柏瑞量化先行混合C"
$detailSheet.Range("H28").Value = 2

# Row 29
$styleSheet.Range("A2").Copy($detailSheet.Range("A29"))
$detailSheet.Range("A29").Value = 27
$detailSheet.Range("B29").NumberFormat = "@"
$detailSheet.Range("B29").Value = "004790"
$detailSheet.Range("B29").Style = "Normal"
$detailSheet.Range("D29").NumberFormat = "@"
$detailSheet.Range("D29").Value = "0.06"
$detailSheet.Range("D29").Style = "Normal"
$detailSheet.Range("E29").NumberFormat = "@"
$detailSheet.Range("E29").Value = "84.00"
$detailSheet.Range("E29").Style = "Normal"
$detailSheet.Range("F29").NumberFormat = "@"
$detailSheet.Range("F29").Value = "1.66"
$detailSheet.Range("F29").Style = "Normal"
$detailSheet.Range("G29").NumberFormat = "@"
$detailSheet.Range("G29").Value = "0.0010"
$detailSheet.Range("G29").Style = "Normal"
$detailSheet.Range("C29").Value = "富荣中证500指数增强A"
$detailSheet.Range("H29").Value = 2

# Row 30
$styleSheet.Range("A2").Copy($detailSheet.Range("A30"))
$detailSheet.Range("A30").Value = 28
$detailSheet.Range("B30").NumberFormat = "@"
$detailSheet.Range("B30").Value = "004791"
$detailSheet.Range("B30").Style = "Normal"
$detailSheet.Range("D30").NumberFormat = "@"
$detailSheet.Range("D30").Value = "0.04"
$detailSheet.Range("D30").Style = "Normal"
$detailSheet.Range("E30").NumberFormat = "@"
$detailSheet.Range("E30").Value = "84.00"
$detailSheet.Range("E30").Style = "Normal"
$detailSheet.Range("F30").NumberFormat = "@"
$detailSheet.Range("F30").Value = "1.66"
$detailSheet.Range("F30").Style = "Normal"
$detailSheet.Range("G30").NumberFormat = "@"
$detailSheet.Range("G30").Value = "0.0007"
$detailSheet.Range("G30").Style = "Normal"
$detailSheet.Range("C30").Value = "富荣中证500指数增强C"
$detailSheet.Range("H30").Value = 2

# ============================================================
# ---- Step 2: append a fresh "总计" sheet with the updated history ----
# ============================================================
$totalSheet = $wb.Worksheets.Add($null, $detailSheet)
$totalSheet.Name = "总计"

$styleSheet.Range("B1").Copy($totalSheet.Range("B1"))
$totalSheet.Range("B1").Value = "日期"
$styleSheet.Range("B1").Copy($totalSheet.Range("C1"))
$totalSheet.Range("C1").Value = "持有数量(只)"
$styleSheet.Range("B1").Copy($totalSheet.Range("D1"))
$totalSheet.Range("D1").Value = "持有市值(亿元)"

# Row 2
$styleSheet.Range("A2").Copy($totalSheet.Range("A2"))
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 29
$totalSheet.Range("D2").Value = 14.7

# Row 3
$styleSheet.Range("A2").Copy($totalSheet.Range("A3"))
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 26
$totalSheet.Range("D3").Value = 3.48

# Row 4
$styleSheet.Range("A2").Copy($totalSheet.Range("A4"))
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 46
$totalSheet.Range("D4").Value = 6.38

# Row 5
$styleSheet.Range("A2").Copy($totalSheet.Range("A5"))
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 74
$totalSheet.Range("D5").Value = 19.9

# Row 6
$styleSheet.Range("A2").Copy($totalSheet.Range("A6"))
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2021-Q1"
$totalSheet.Range("C6").Value = 41
$totalSheet.Range("D6").Value = 16.19

# Row 7
$styleSheet.Range("A2").Copy($totalSheet.Range("A7"))
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("B7").Value = "2020-Q4"
$totalSheet.Range("C7").Value = 14
$totalSheet.Range("D7").Value = 6.6
